# Auto-generated: reproduces the scheduled-runner market-data refresh diff
# for Sheets/Bahamut_Profits.xlsx (workbook tab layout: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2229.6
$ws.Range("I6").Value = 2537
$ws.Range("K6").Value = 7611
$ws.Range("M6").Value = -7499
$ws.Range("H19").Value = 1670.6957
$ws.Range("I19").Value = 1438.4
$ws.Range("J19").Value = 1849.3846
$ws.Range("K19").Value = 1438.4
$ws.Range("L19").Value = 1849.3846
$ws.Range("M19").Value = -1263.4
$ws.Range("N19").Value = -2199.3846
$ws.Range("H33").Value = 114.25
$ws.Range("I33").Value = 107.181816
$ws.Range("K33").Value = 107.181816
$ws.Range("M33").Value = 121.818184
$ws.Range("H40").Value = 26317990
$ws.Range("I40").Value = 1910
$ws.Range("J40").Value = 35716590
$ws.Range("K40").Value = 1910
$ws.Range("L40").Value = 35716590
$ws.Range("M40").Value = -1735
$ws.Range("N40").Value = -35716940
$ws.Range("H112").Value = 1378.0769
$ws.Range("I112").Value = 1237.5
$ws.Range("J112").Value = 1403.6364
$ws.Range("K112").Value = 3712.5
$ws.Range("L112").Value = 4210.9092
$ws.Range("M112").Value = -2604.5
$ws.Range("N112").Value = -6426.9092
$ws.Range("H125").Value = 3259.4211
$ws.Range("I125").Value = 1753.3077
$ws.Range("K125").Value = 15779.7693
$ws.Range("M125").Value = -13319.7693
$ws.Range("H138").Value = 2917.74
$ws.Range("I138").Value = 667.7037
$ws.Range("J138").Value = 3749.9453
$ws.Range("K138").Value = 2003.1111
$ws.Range("L138").Value = 11249.8359
$ws.Range("N138").Value = -21529.8359
$ws.Range("M138").Value = 3136.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2742.6316
$ws.Range("I2").Value = 2497.4138
$ws.Range("J2").Value = 3532.7778
$ws.Range("K2").Value = 2497.4138
$ws.Range("L2").Value = 3532.7778
$ws.Range("M2").Value = -2384.4138
$ws.Range("N2").Value = -3758.7778
$ws.Range("H61").Value = 1856.3334
$ws.Range("I61").Value = 1808.4445
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1808.4445
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1596.4445
$ws.Range("N61").Value = -2424
$ws.Range("H116").Value = 2742.6316
$ws.Range("I116").Value = 2497.4138
$ws.Range("J116").Value = 3532.7778
$ws.Range("K116").Value = 2497.4138
$ws.Range("L116").Value = 3532.7778
$ws.Range("M116").Value = -203.4137999999998
$ws.Range("N116").Value = -8120.7778
$ws.Range("H132").Value = 1726
$ws.Range("I132").Value = 1180.421
$ws.Range("J132").Value = 3799.2
$ws.Range("K132").Value = 3541.263
$ws.Range("L132").Value = 11397.6
$ws.Range("M132").Value = -1011.263
$ws.Range("N132").Value = -16457.6
$ws.Range("H136").Value = 1856.3334
$ws.Range("I136").Value = 1808.4445
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5425.333500000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2875.333500000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2742.6316
$ws.Range("I3").Value = 2497.4138
$ws.Range("J3").Value = 3532.7778
$ws.Range("K3").Value = 2497.4138
$ws.Range("L3").Value = 3532.7778
$ws.Range("M3").Value = -2383.4138
$ws.Range("N3").Value = -3760.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 234353.44
$ws.Range("I99").Value = 346389.97
$ws.Range("J99").Value = 2277.7856
$ws.Range("K99").Value = 346389.97
$ws.Range("L99").Value = 2277.7856
$ws.Range("M99").Value = -344891.97
$ws.Range("N99").Value = -5273.7856
$ws.Range("H126").Value = 234353.44
$ws.Range("I126").Value = 346389.97
$ws.Range("J126").Value = 2277.7856
$ws.Range("K126").Value = 1039169.91
$ws.Range("L126").Value = 6833.3568
$ws.Range("M126").Value = -1036699.91
$ws.Range("N126").Value = -11773.3568
$ws.Range("H132").Value = 2756
$ws.Range("I132").Value = 2155.7827
$ws.Range("J132").Value = 4289.8887
$ws.Range("K132").Value = 6467.348100000001
$ws.Range("L132").Value = 12869.6661
$ws.Range("M132").Value = -3937.348100000001
$ws.Range("N132").Value = -17929.6661
$ws.Range("H134").Value = 2638.9312
$ws.Range("I134").Value = 2356.7778
$ws.Range("J134").Value = 3100.6365
$ws.Range("K134").Value = 7070.3334
$ws.Range("L134").Value = 9301.9095
$ws.Range("M134").Value = -4535.3334
$ws.Range("N134").Value = -14371.9095
$ws.Range("H140").Value = 55406
$ws.Range("J140").Value = 55406
$ws.Range("L140").Value = 55406
$ws.Range("N140").Value = -65766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 2964.6667
$ws.Range("I56").Value = 2964.6667
$ws.Range("K56").Value = 2964.6667
$ws.Range("M56").Value = -2434.6667
$ws.Range("H64").Value = 1594.9231
$ws.Range("I64").Value = 964
$ws.Range("J64").Value = 1784.2
$ws.Range("K64").Value = 2892
$ws.Range("L64").Value = 5352.6
$ws.Range("M64").Value = -2622
$ws.Range("N64").Value = -5892.6
$ws.Range("H67").Value = 1594.9231
$ws.Range("I67").Value = 964
$ws.Range("J67").Value = 1784.2
$ws.Range("K67").Value = 2892
$ws.Range("L67").Value = 5352.6
$ws.Range("M67").Value = -1956
$ws.Range("N67").Value = -7224.6
$ws.Range("H87").Value = 12484
$ws.Range("I87").Value = 3393
$ws.Range("J87").Value = 21575
$ws.Range("K87").Value = 10179
$ws.Range("L87").Value = 64725
$ws.Range("M87").Value = -8931
$ws.Range("N87").Value = -67221
$ws.Range("H90").Value = 12484
$ws.Range("I90").Value = 3393
$ws.Range("J90").Value = 21575
$ws.Range("K90").Value = 30537
$ws.Range("L90").Value = 194175
$ws.Range("M90").Value = -24297
$ws.Range("N90").Value = -206655
$ws.Range("H117").Value = 2552.5293
$ws.Range("J117").Value = 2857.5334
$ws.Range("L117").Value = 8572.600199999999
$ws.Range("N117").Value = -15456.6002
$ws.Range("H124").Value = 4375
$ws.Range("I124").Value = 1250
$ws.Range("K124").Value = 3750
$ws.Range("M124").Value = 1160
$ws.Range("H129").Value = 2616.7368
$ws.Range("I129").Value = 1298.1666
$ws.Range("J129").Value = 3225.3076
$ws.Range("K129").Value = 3894.4998
$ws.Range("L129").Value = 9675.9228
$ws.Range("M129").Value = 1105.5002
$ws.Range("N129").Value = -19675.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1367.7084
$ws.Range("I102").Value = 1084.7894
$ws.Range("J102").Value = 2442.8
$ws.Range("K102").Value = 1084.7894
$ws.Range("L102").Value = 2442.8
$ws.Range("M102").Value = 537.2106000000001
$ws.Range("N102").Value = -5686.8
$ws.Range("H126").Value = 2992.4
$ws.Range("I126").Value = 2964.3157
$ws.Range("K126").Value = 8892.947100000001
$ws.Range("M126").Value = -6422.947100000001
$ws.Range("H132").Value = 2913.75
$ws.Range("I132").Value = 956
$ws.Range("K132").Value = 2868
$ws.Range("M132").Value = -338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3949.7585
$ws.Range("I136").Value = 1301.9
$ws.Range("J136").Value = 9833.888999999999
$ws.Range("K136").Value = 3905.7
$ws.Range("L136").Value = 29501.667
$ws.Range("M136").Value = -1355.7
$ws.Range("N136").Value = -34601.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9557.333000000001
$ws.Range("J2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("N2").Value = -6224
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H126").Value = 502.73077
$ws.Range("I126").Value = 459
$ws.Range("J126").Value = 601.125
$ws.Range("K126").Value = 1377
$ws.Range("L126").Value = 1803.375
$ws.Range("M126").Value = 1093
$ws.Range("N126").Value = -6743.375
$ws.Range("H138").Value = 53700
$ws.Range("J138").Value = 53700
$ws.Range("L138").Value = 53700
$ws.Range("N138").Value = -63980
